$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 ("Rules" sheet) changes from the text "R40" to the text "1".
# Build the new value in an unused scratch cell via a formula so it is
# produced as TEXT (not a number), then copy/paste-values onto B11 so the
# cell's existing style (s="23") and number format are left untouched.
$scratch = $ws.Range("Z1")
$scratch.Formula = '=TEXT(1,"0")'
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
